# Apply the RoguishArchetypes.xlsx edit described by the diff:
#  - Row 1 keeps only A1/B1 (Arcane Trickster / Assassin); C1:K1 removed.
#  - New rows 2-8 added in columns A and B with the Arcane Trickster archetype data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused header cells C1:K1 entirely (not just their contents)
$ws.Range("C1:K1").Clear()

# Make sure the new cells A2:B8 inherit the same style (s="1") as A1/B1
$ws.Range("A1").Copy()
$ws.Range("A2:B8").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = "None"
$ws.Range("B2").Value = "None"

# Row 3
$ws.Range("A3").Value = "None"
$ws.Range("B3").Value = "None"

# Row 4
$ws.Range("A4").Value = "None"
$ws.Range("B4").Value = "Disguise Kit=Poisoner's Kit"

# Row 5
$ws.Range("A5").Value = "None"
$ws.Range("B5").Value = "None"

# Row 6
$ws.Range("A6").Value = "true=classes/rogue/BonusSpellsArcaneTrickster.xlsx"
$ws.Range("B6").Value = "false="

# Row 7
$ws.Range("A7").Value = "true=Intelligence=classes/rogue/SpellSlotsArcaneTrickster.xlsx=classes/rogue/SpellListArcaneTrickster.xlsx=K"
$ws.Range("B7").Value = "false="

# Row 8
$ws.Range("A8").Value = "3/Mage Hand=3/Mage Hand Legerdemain=9/Magical Ambush=13/Versatile Trickster=17/Spell Thief"
$ws.Range("B8").Value = "3/Assassinate=9/Infiltration Expertise=13/Impostor=17/Death Strike"
